$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 709.5
$ws.Range("I135").Value = 709.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6385.5
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = -3850.5
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5223.6665
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5223.6665
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = 5223.6665
$ws.Range("N61").Value = -5647.6665
$ws.Range("L61").ClearContents()
$ws.Range("H121").Value = 29681.818
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 29681.818
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 29681.818
$ws.Range("N121").Value = -33175.818
$ws.Range("H122").Value = 1196.3793
$ws.Range("I122").Value = 967.9545000000001
$ws.Range("J122").Value = 1914.2858
$ws.Range("K122").Value = 2903.8635
$ws.Range("L122").Value = 5742.857400000001
$ws.Range("M122").Value = -453.8635000000004
$ws.Range("N122").Value = -10642.8574
$ws.Range("H123").Value = 32402.857
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 32402.857
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 32402.857
$ws.Range("N123").Value = -42202.857
$ws.Range("H124").Value = 18729.691
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 18729.691
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 18729.691
$ws.Range("N124").Value = -28549.691
$ws.Range("H125").Value = 31244.75
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 31244.75
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 31244.75
$ws.Range("N125").Value = -41084.75
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 34000
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 34000
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 34000
$ws.Range("N127").Value = -43920
$ws.Range("H128").Value = 34280
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 34280
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 34280
$ws.Range("N128").Value = -44240
$ws.Range("H129").Value = 39333
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 39333
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 39333
$ws.Range("N129").Value = -49333
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 1882.0209
$ws.Range("I132").Value = 1474.875
$ws.Range("J132").Value = 3917.75
$ws.Range("K132").Value = 4424.625
$ws.Range("L132").Value = 11753.25
$ws.Range("M132").Value = -1894.625
$ws.Range("N132").Value = -16813.25
$ws.Range("H133").Value = 100000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 100000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 28016.666
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 28016.666
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 28016.666
$ws.Range("N135").Value = -38156.666
$ws.Range("H136").Value = 5223.6665
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 5223.6665
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = 15670.9995
$ws.Range("N136").Value = -20770.9995
$ws.Range("L136").ClearContents()
$ws.Range("H137").Value = 59166.668
$ws.Range("I137").Value = 39000
$ws.Range("J137").Value = 63200
$ws.Range("K137").Value = 39000
$ws.Range("L137").Value = 63200
$ws.Range("M137").Value = -33900
$ws.Range("N137").Value = -73400
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 300000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 300000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 300000
$ws.Range("N139").Value = -310280
$ws.Range("H140").Value = 35325
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 35325
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 35325
$ws.Range("N140").Value = -45685
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5035.048
$ws.Range("I134").Value = 5362.4
$ws.Range("J134").Value = 4216.6665
$ws.Range("K134").Value = 16087.2
$ws.Range("L134").Value = 12649.9995
$ws.Range("M134").Value = -13552.2
$ws.Range("N134").Value = -17719.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13744.25
$ws.Range("I58").Value = 17825.666
$ws.Range("J58").Value = 1500
$ws.Range("K58").Value = 17825.666
$ws.Range("L58").Value = 1500
$ws.Range("M58").Value = -17622.666
$ws.Range("N58").Value = -1906
$ws.Range("H132").Value = 2705.8667
$ws.Range("I132").Value = 2519.9688
$ws.Range("J132").Value = 3163.4614
$ws.Range("K132").Value = 7559.9064
$ws.Range("L132").Value = 9490.3842
$ws.Range("M132").Value = -5029.9064
$ws.Range("N132").Value = -14550.3842
$ws.Range("H136").Value = 13744.25
$ws.Range("I136").Value = 17825.666
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 53476.99800000001
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -50926.99800000001
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3033.3333
$ws.Range("I61").Value = 2500
$ws.Range("J61").Value = 3300
$ws.Range("K61").Value = 2500
$ws.Range("L61").Value = 3300
$ws.Range("M61").Value = -2298
$ws.Range("N61").Value = -3704
$ws.Range("H113").Value = 3033.3333
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 3300
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3300
$ws.Range("M113").Value = -330
$ws.Range("N113").Value = -7640
$ws.Range("H132").Value = 4217.7837
$ws.Range("I132").Value = 3747.7273
$ws.Range("J132").Value = 4907.2
$ws.Range("K132").Value = 11243.1819
$ws.Range("L132").Value = 14721.6
$ws.Range("M132").Value = -8713.1819
$ws.Range("N132").Value = -19781.6
